$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '72.420.13'
$ws.Range('E2').Value = '  +1.73%  '
$ws.Range('D3').Value = '2.687.10'
$ws.Range('E3').Value = '  +2.11%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '601.65'
$ws.Range('E5').Value = '  -0.83%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '177.76'
$ws.Range('E6').Value = '  -1.74%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').Value = '2.684.98'
$ws.Range('E9').Value = '  +2.16%  '
$ws.Range('E10').Value = '  +3.03%  '
$ws.Range('E11').Value = '  +2.41%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.355'
$ws.Range('E12').Value = '  +2.44%  '
$ws.Range('E13').Value = '  +0.43%  '
$ws.Range('D14').Value = '3.161.55'
$ws.Range('E14').Value = '  +2.82%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000186'
$ws.Range('E15').Value = '  -0.78%  '
$ws.Range('D16').Value = '72.324.33'
$ws.Range('E16').Value = '  +1.71%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.40'
$ws.Range('E17').Value = '  -1.11%  '
$ws.Range('D18').Value = '2.686.87'
$ws.Range('E18').Value = '  +2.44%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.06'
$ws.Range('E19').Value = '  +4.72%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.94'
$ws.Range('E20').Value = '  +0.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '371.69'
$ws.Range('E21').Value = '  -3.19%  '
$ws.Range('E22').Value = '  +0.93%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.06'
$ws.Range('E23').Value = '  +8.34%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '72.27'
$ws.Range('E24').Value = '  +0.16%  '
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('E26').Value = '  -1.99%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.85'
$ws.Range('E27').Value = '  +2.06%  '
$ws.Range('D28').Value = '2.827.35'
$ws.Range('E28').Value = '  +2.36%  '
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('D30').Value = '0.0₃0944'
$ws.Range('E30').Value = '  -2.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.11'
$ws.Range('E31').Value = '  +0.63%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '513.48'
$ws.Range('E32').Value = '  -5.24%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.30'
$ws.Range('E33').Value = '  -1.80%  '
$ws.Range('E34').Value = '  -0.93%  '
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '164.40'
$ws.Range('E36').Value = '  -0.10%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.62'
$ws.Range('E37').Value = '  +1.94%  '
$ws.Range('E38').Value = '  +0.50%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.38'
$ws.Range('E39').Value = '  -0.17%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.82'
$ws.Range('E40').Value = '  -3.55%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.108'
$ws.Range('E41').Value = '  -8.81%  '
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.03'
$ws.Range('E43').Value = '  -0.39%  '
$ws.Range('E44').Value = '  -2.59%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.333'
$ws.Range('E45').Value = '  +0.71%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '39.27'
$ws.Range('E46').Value = '  -1.64%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '153.75'
$ws.Range('E47').Value = '  -0.31%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.74'
$ws.Range('E48').Value = '  +2.68%  '
$ws.Range('E49').Value = '  +3.77%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.73'
$ws.Range('E50').Value = '  +1.91%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0767'
$ws.Range('E51').Value = '  +1.77%  '
